# Auto-generated Excel COM-interop script
# Implements S27/G01: Portfolio rebalance v1: target-weight + budget
# (plus the accompanying S26 row formatting/cleanup changes)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 247-254: row height tweaks + style/cleanup on existing cells ---
$ws.Rows.Item(247).RowHeight = 41.25
$ws.Rows.Item(254).RowHeight = 27.75

$ws.Rows.Item(248).RowHeight = 41.25
$c = $ws.Range("H248")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$ws.Range("I248").Clear()

$ws.Rows.Item(249).RowHeight = 41.25
$c = $ws.Range("H249")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$ws.Range("I249").Clear()

$ws.Rows.Item(250).RowHeight = 41.25
$c = $ws.Range("F250")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("H250")
$c.WrapText = $true
$c.VerticalAlignment = -4160

$ws.Rows.Item(251).RowHeight = 41.25
$c = $ws.Range("H251")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$ws.Range("I251").Clear()

$ws.Rows.Item(252).RowHeight = 41.25
$c = $ws.Range("H252")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$ws.Range("I252").Clear()

$ws.Rows.Item(253).RowHeight = 41.25
$c = $ws.Range("H253")
$c.WrapText = $true
$c.VerticalAlignment = -4160

# --- Rows 255-260: apply wrap/top styling + row heights (content unchanged) ---
$rng = $ws.Range("A255:H255")
$rng.WrapText = $true
$rng.VerticalAlignment = -4160
$ws.Rows.Item(255).RowHeight = 41.75

$rng = $ws.Range("A256:H256")
$rng.WrapText = $true
$rng.VerticalAlignment = -4160
$ws.Rows.Item(256).RowHeight = 41.75

$rng = $ws.Range("A257:H257")
$rng.WrapText = $true
$rng.VerticalAlignment = -4160
$ws.Rows.Item(257).RowHeight = 28.35

$rng = $ws.Range("A258:H258")
$rng.WrapText = $true
$rng.VerticalAlignment = -4160
$ws.Rows.Item(258).RowHeight = 28.35

$rng = $ws.Range("A259:H259")
$rng.WrapText = $true
$rng.VerticalAlignment = -4160
$ws.Rows.Item(259).RowHeight = 41.75

$rng = $ws.Range("A260:H260")
$rng.WrapText = $true
$rng.VerticalAlignment = -4160
$ws.Rows.Item(260).RowHeight = 28.35

# --- Rows 261-280: new S27 rows (Portfolio rebalance v1/v2/v3 epics) ---
# Row 261
$c = $ws.Range("A261")
$c.Value = ("S27")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("B261")
$c.Value = ("G01")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("C261")
$c.Value = ("Portfolio rebalance v1: target-weight + budget")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("D261")
$c.Value = ("S27_G01_TB001")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("E261")
$c.Value = ("Backend: Add DB schema for portfolio rebalance (policy + schedule + run history + run orders) scoped to user and group.")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("F261")
$c.Value = ("See docs/portfolio_rebalancing.md. Keep policy snapshot per run for audit.")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("G261")
$c.Value = ("implemented")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("H261")
$c.Value = ("Rebalance schema added (policy/schedule/run/run_orders) + Alembic migration 0040.")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$ws.Rows.Item(261).RowHeight = 41.75

# Row 262
$c = $ws.Range("A262")
$c.Value = ("S27")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("B262")
$c.Value = ("G01")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("C262")
$c.Value = ("Portfolio rebalance v1: target-weight + budget")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("D262")
$c.Value = ("S27_G01_TB002")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("E262")
$c.Value = ("Backend: Implement rebalance computation engine (live weights, drift bands, budget scaling, qty rounding, min trade value, max trades).")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("F262")
$c.Value = ("Supports broker-scoped runs (Zerodha/AngelOne) and " + [char]0x201C + "Both (run separately)" + [char]0x201D + ".")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("G262")
$c.Value = ("implemented")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("H262")
$c.Value = ("Deterministic target-weight rebalance engine with drift bands, budget scaling, qty rounding, min trade value, max trades, and before/after drift metrics.")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$ws.Rows.Item(262).RowHeight = 41.75

# Row 263
$c = $ws.Range("A263")
$c.Value = ("S27")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("B263")
$c.Value = ("G01")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("C263")
$c.Value = ("Portfolio rebalance v1: target-weight + budget")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("D263")
$c.Value = ("S27_G01_TB003")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("E263")
$c.Value = ("Backend: Add /api/rebalance/preview endpoint returning proposed trades + summary metrics (turnover, drift reduced, budget used).")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("F263")
$c.Value = ("Input: group_id, broker_name, budget_pct/amount, bands, constraints.")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("G263")
$c.Value = ("implemented")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("H263")
$c.Value = ("Added /api/rebalance/preview returning per-broker proposals + summary + warnings.")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$ws.Rows.Item(263).RowHeight = 41.75

# Row 264
$c = $ws.Range("A264")
$c.Value = ("S27")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("B264")
$c.Value = ("G01")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("C264")
$c.Value = ("Portfolio rebalance v1: target-weight + budget")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("D264")
$c.Value = ("S27_G01_TB004")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("E264")
$c.Value = ("Backend: Add /api/rebalance/execute endpoint to create orders (manual queue or AUTO) and persist run history + order mapping.")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("F264")
$c.Value = ("Idempotency key recommended; store execution results/status per order.")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("G264")
$c.Value = ("implemented")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("H264")
$c.Value = ("Added /api/rebalance/execute to create orders (MANUAL/AUTO) and persist RebalanceRun + RebalanceRunOrder mappings; idempotency supported.")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$ws.Rows.Item(264).RowHeight = 41.75

# Row 265
$c = $ws.Range("A265")
$c.Value = ("S27")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("B265")
$c.Value = ("G01")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("C265")
$c.Value = ("Portfolio rebalance v1: target-weight + budget")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("D265")
$c.Value = ("S27_G01_TF001")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("E265")
$c.Value = ("Frontend: Add Rebalance panel on portfolio group views (next rebalance, last rebalance, broker scope, preview/run/history).")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("F265")
$c.Value = ("Surface core knobs: budget %, drift bands, max trades, min trade value.")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("G265")
$c.Value = ("implemented")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("H265")
$c.Value = ("Added Rebalance action for PORTFOLIO groups in Groups page with preview/run/history workflow.")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$ws.Rows.Item(265).RowHeight = 41.75

# Row 266
$c = $ws.Range("A266")
$c.Value = ("S27")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("B266")
$c.Value = ("G01")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("C266")
$c.Value = ("Portfolio rebalance v1: target-weight + budget")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("D266")
$c.Value = ("S27_G01_TF002")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("E266")
$c.Value = ("Frontend: Rebalance preview dialog (proposed buys/sells, before/after weights, budget scaling, warnings, confirm).")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("F266")
$c.Value = ("Reuse existing order UX patterns; confirm can route to manual queue or AUTO broker.")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("G266")
$c.Value = ("implemented")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("H266")
$c.Value = ("Preview dialog shows proposed buys/sells, before/after drift summary, warnings, and confirm actions (manual queue or AUTO).")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$ws.Rows.Item(266).RowHeight = 41.75

# Row 267
$c = $ws.Range("A267")
$c.Value = ("S27")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("B267")
$c.Value = ("G01")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("C267")
$c.Value = ("Portfolio rebalance v1: target-weight + budget")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("D267")
$c.Value = ("S27_G01_TF003")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("E267")
$c.Value = ("Frontend: Rebalance history tab/table + run detail (inputs snapshot, proposed orders, execution status/errors).")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("F267")
$c.Value = ("Link to Orders/Queue filtered by rebalance run id.")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("G267")
$c.Value = ("implemented")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("H267")
$c.Value = ("History tab lists recent runs and run-order details (order ids, symbols, qty, status).")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("I267")
$c.Value = ("Optional: add Orders/Queue filtering by rebalance run id.")
$ws.Rows.Item(267).RowHeight = 41.75

# Row 268
$c = $ws.Range("A268")
$c.Value = ("S27")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("B268")
$c.Value = ("G01")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("C268")
$c.Value = ("Portfolio rebalance v1: target-weight + budget")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("D268")
$c.Value = ("S27_G01_TT001")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("E268")
$c.Value = ("Tests: Add backend tests for preview/execute endpoints (budget scaling, band gating, rounding, max trades, broker-scope).")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("G268")
$c.Value = ("implemented")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("H268")
$c.Value = ("Added backend tests for preview/execute (budget scaling + idempotency) with broker calls monkeypatched.")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$ws.Rows.Item(268).RowHeight = 41.75

# Row 269
$c = $ws.Range("A269")
$c.Value = ("S27")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("B269")
$c.Value = ("G02")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("C269")
$c.Value = ("Rebalance scheduling + next opportunity")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("D269")
$c.Value = ("S27_G02_TB001")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("E269")
$c.Value = ("Backend: Implement schedule model + next_rebalance_at computation (weekly/monthly/quarterly/custom days, timezone Asia/Kolkata).")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("F269")
$c.Value = ("No auto-execution in v1; just compute/display next opportunity.")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("G269")
$c.Value = ("planned")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("H269")
$c.Value = ("Schedule stored per portfolio group; next opportunity computed consistently.")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$ws.Rows.Item(269).RowHeight = 55.2

# Row 270
$c = $ws.Range("A270")
$c.Value = ("S27")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("B270")
$c.Value = ("G02")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("C270")
$c.Value = ("Rebalance scheduling + next opportunity")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("D270")
$c.Value = ("S27_G02_TF001")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("E270")
$c.Value = ("Frontend: Schedule editor in Rebalance panel (frequency, time, timezone) + display next/last rebalance timestamps.")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("G270")
$c.Value = ("planned")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("H270")
$c.Value = ("User can configure schedule and see next rebalance date/opportunity.")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$ws.Rows.Item(270).RowHeight = 41.75

# Row 271
$c = $ws.Range("A271")
$c.Value = ("S27")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("B271")
$c.Value = ("G02")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("C271")
$c.Value = ("Rebalance scheduling + next opportunity")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("D271")
$c.Value = ("S27_G02_TB002")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("E271")
$c.Value = ("Backend: Expose schedule/policy read+update endpoints for portfolio groups (GET/PUT).")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("G271")
$c.Value = ("planned")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("H271")
$c.Value = ("FE can persist schedule/policy config per group.")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$ws.Rows.Item(271).RowHeight = 28.35

# Row 272
$c = $ws.Range("A272")
$c.Value = ("S27")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("B272")
$c.Value = ("G02")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("C272")
$c.Value = ("Rebalance scheduling + next opportunity")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("D272")
$c.Value = ("S27_G02_TB003")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("E272")
$c.Value = ("Backend (optional): Add background scheduler (disabled by default) to auto-run rebalance previews/notifications and/or execute when enabled.")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("F272")
$c.Value = ("Defer if needed; keep behind settings flag.")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("G272")
$c.Value = ("planned")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("H272")
$c.Value = ("Foundation for automated periodic rebalancing.")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$ws.Rows.Item(272).RowHeight = 41.75

# Row 273
$c = $ws.Range("A273")
$c.Value = ("S27")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("B273")
$c.Value = ("G03")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("C273")
$c.Value = ("Rebalance v2: signal/strategy-driven rotation")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("D273")
$c.Value = ("S27_G03_TB001")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("E273")
$c.Value = ("Backend: Extend rebalance engine to support signal-driven targets (use saved Strategy output to rank universe and derive weights).")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("F273")
$c.Value = ("Inputs: strategy_id/version/output, candidate universe group/screener run, top-N, weighting scheme.")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("G273")
$c.Value = ("planned")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("H273")
$c.Value = ("Preview returns target weights derived from strategy + resulting trades.")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$ws.Rows.Item(273).RowHeight = 41.75

# Row 274
$c = $ws.Range("A274")
$c.Value = ("S27")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("B274")
$c.Value = ("G03")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("C274")
$c.Value = ("Rebalance v2: signal/strategy-driven rotation")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("D274")
$c.Value = ("S27_G03_TF001")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("E274")
$c.Value = ("Frontend: Add " + [char]0x201C + "Signal-driven rebalance" + [char]0x201D + " mode UI (select strategy/version/output, candidate universe, top-N, weighting scheme).")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("G274")
$c.Value = ("planned")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("H274")
$c.Value = ("User can rotate portfolio using strategy signals with preview/run flow.")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$ws.Rows.Item(274).RowHeight = 41.75

# Row 275
$c = $ws.Range("A275")
$c.Value = ("S27")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("B275")
$c.Value = ("G03")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("C275")
$c.Value = ("Rebalance v2: signal/strategy-driven rotation")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("D275")
$c.Value = ("S27_G03_TB002")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("E275")
$c.Value = ("Backend: Add replacement rules (sell criteria, eligibility filters, min liquidity/price, blacklist/whitelist) and audit reasons.")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("F275")
$c.Value = ("Keep knobs configurable per policy.")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("G275")
$c.Value = ("planned")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("H275")
$c.Value = ("Rotation decisions are explainable and consistent.")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$ws.Rows.Item(275).RowHeight = 41.75

# Row 276
$c = $ws.Range("A276")
$c.Value = ("S27")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("B276")
$c.Value = ("G03")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("C276")
$c.Value = ("Rebalance v2: signal/strategy-driven rotation")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("D276")
$c.Value = ("S27_G03_TT001")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("E276")
$c.Value = ("Tests: Coverage for strategy-driven target derivation + replacement rules and constraints.")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("G276")
$c.Value = ("planned")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("H276")
$c.Value = ("Ensure deterministic outputs given fixed candle inputs.")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$ws.Rows.Item(276).RowHeight = 28.35

# Row 277
$c = $ws.Range("A277")
$c.Value = ("S27")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("B277")
$c.Value = ("G04")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("C277")
$c.Value = ("Rebalance v3: risk-based (risk parity / contributions)")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("D277")
$c.Value = ("S27_G04_TB001")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("E277")
$c.Value = ("Backend: Compute/caches covariance matrix from price history and per-asset risk metrics needed for risk contribution.")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("F277")
$c.Value = ("Use 6M" + [char]0x2013 + "1Y windows; cache keyed by date/window/universe.")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("G277")
$c.Value = ("planned")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("H277")
$c.Value = ("Covariance inputs available for optimizer.")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$ws.Rows.Item(277).RowHeight = 41.75

# Row 278
$c = $ws.Range("A278")
$c.Value = ("S27")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("B278")
$c.Value = ("G04")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("C278")
$c.Value = ("Rebalance v3: risk-based (risk parity / contributions)")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("D278")
$c.Value = ("S27_G04_TB002")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("E278")
$c.Value = ("Backend: Implement risk-based optimizer (risk parity / equal risk contribution) with constraints (max weight, max trades, budget).")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("F278")
$c.Value = ("Start with a simple iterative solver; avoid heavy deps in v1.")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("G278")
$c.Value = ("planned")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("H278")
$c.Value = ("Target weights produced from risk objective; preview returns resulting trades.")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$ws.Rows.Item(278).RowHeight = 41.75

# Row 279
$c = $ws.Range("A279")
$c.Value = ("S27")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("B279")
$c.Value = ("G04")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("C279")
$c.Value = ("Rebalance v3: risk-based (risk parity / contributions)")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("D279")
$c.Value = ("S27_G04_TF001")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("E279")
$c.Value = ("Frontend: Add " + [char]0x201C + "Risk-based rebalance" + [char]0x201D + " mode UI (window, constraints, objective) within the same preview/run workflow.")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("G279")
$c.Value = ("planned")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("H279")
$c.Value = ("Risk-based rebalance is configurable and consistent with other modes.")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$ws.Rows.Item(279).RowHeight = 41.75

# Row 280
$c = $ws.Range("A280")
$c.Value = ("S27")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("B280")
$c.Value = ("G04")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("C280")
$c.Value = ("Rebalance v3: risk-based (risk parity / contributions)")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("D280")
$c.Value = ("S27_G04_TT001")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("E280")
$c.Value = ("Tests: Add unit tests for covariance estimation + optimizer sanity checks + performance guardrails.")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("G280")
$c.Value = ("planned")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$c = $ws.Range("H280")
$c.Value = ("Confident correctness and reasonable runtime.")
$c.WrapText = $true
$c.VerticalAlignment = -4160
$ws.Rows.Item(280).RowHeight = 28.35

# --- View state: scroll position + active cell (best effort) ---
$win = $excel.ActiveWindow
$win.ScrollRow = 248
$win.ScrollColumn = 1
$ws.Range("E259").Select() | Out-Null

